$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.254.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.598.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'212.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.244"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.0606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'18.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "'0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "'1.824.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "'1.597.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'63.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "'26.253.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'227.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.97%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'8.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'145.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'6.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'15.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'1.444.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "'0.820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'5.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").Value = "'0.930"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").Value = "'1.736.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'0.754"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "'60.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'87.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "'1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'0.0499"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0948"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.16%  "
